# BUG: Fix read_excel w/parse_cols & empty dataset (pandas gh-9208)
#
# Adds a third worksheet ("Sheet3") to the workbook that has a header row
# (A, B, C, D, E, F) but no data rows underneath it -- this is the
# "empty dataset" fixture used to exercise read_excel(parse_cols=...)
# against a sheet that only contains headers.

$wb = $excel.ActiveWorkbook

# Append the new sheet after the existing ones (Sheet1, Sheet2) so it
# becomes the last/active tab, matching the target workbook layout.
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws3 = $wb.Worksheets.Add($null, $lastSheet)
$ws3.Name = "Sheet3"

# Header-only row; no data rows follow.
$ws3.Range("A1").Value = "A"
$ws3.Range("B1").Value = "B"
$ws3.Range("C1").Value = "C"
$ws3.Range("D1").Value = "D"
$ws3.Range("E1").Value = "E"
$ws3.Range("F1").Value = "F"

# Leave the cursor one cell past the typed header row (F2), matching the
# recorded selection state of the source workbook.
$ws3.Range("F2").Select()
